# Update player roster table (A2:C18) with refreshed 2024-25 season data.
# Columns: A = Oyuncu Adı (Player), B = Pozisyon (Position), C = Takım (Team)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Derrick White",      "PG,SG",    "Boston Celtics"),
    @("Cade Cunningham",    "PG,SG",    "Detroit Pistons"),
    @("LaMelo Ball",        "PG,SG",    "Charlotte Hornets"),
    @("Isaiah Collier",     "PG",       "Utah Jazz"),
    @("Julius Randle",      "PF,C",     "Minnesota Timberwolves"),
    @("Naz Reid",           "PF,C",     "Minnesota Timberwolves"),
    @("Onyeka Okongwu",     "PF,C",     "Atlanta Hawks"),
    @("Cole Anthony",       "PG",       "Orlando Magic"),
    @("Damian Lillard",     "PG",       "Milwaukee Bucks"),
    @("Deandre Ayton",      "C",        "Portland Trail Blazers"),
    @("Bam Adebayo",        "C",        "Miami Heat"),
    @("Jerami Grant",       "SF,PF",    "Portland Trail Blazers"),
    @("Malik Monk",         "PG,SG,SF", "Sacramento Kings"),
    @("Cameron Johnson",    "SF,PF",    "Brooklyn Nets"),
    @("Anthony Davis",      "PF,C",     "Los Angeles Lakers"),
    @("Isaiah Hartenstein",  "C",        "Oklahoma City Thunder"),
    @("Collin Sexton",      "PG,SG",    "Utah Jazz")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
